$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.531.95"
$ws.Range("E2").Value = "  +0.93%  "

$ws.Range("D3").Value = "3.385.93"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'575.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "

$ws.Range("D6").Value = "'140.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.56%  "

$ws.Range("E8").Value = "  -0.55%  "

$ws.Range("D9").Value = "'7.69"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.74%  "

$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("D11").Value = "'0.385"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.64%  "

$ws.Range("D12").Value = "3.966.04"
$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("E13").Value = "  +0.17%  "

$ws.Range("D14").Value = "'28.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.37%  "

$ws.Range("D15").Value = "3.392.62"
$ws.Range("E15").Value = "  -0.12%  "

$ws.Range("E16").Value = "  -0.80%  "

$ws.Range("D17").Value = "61.556.05"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("D18").Value = "'6.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.51%  "

$ws.Range("D19").Value = "'13.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.87%  "

$ws.Range("D20").Value = "'8.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("D21").Value = "'390.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.76%  "

$ws.Range("D22").Value = "'75.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.25%  "

$ws.Range("D23").Value = "'0.552"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.36%  "

$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("E25").Value = "  -5.22%  "

$ws.Range("E26").Value = "  +7.02%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").Value = "'7.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.07%  "

$ws.Range("E29").Value = "  +0.28%  "

$ws.Range("E30").Value = "  -0.89%  "

$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("D32").Value = "'1.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.82%  "

$ws.Range("D33").Value = "'23.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.36%  "

$ws.Range("D34").Value = "'6.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.46%  "

$ws.Range("D35").Value = "'168.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.46%  "

$ws.Range("D36").Value = "'5.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.55%  "

$ws.Range("D37").Value = "3.420.55"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").Value = "'1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.68%  "

$ws.Range("D39").Value = "'0.0766"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.83%  "

$ws.Range("D40").Value = "'26.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.91%  "

$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("D42").Value = "'4.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.59%  "

$ws.Range("E43").Value = "  -1.36%  "

$ws.Range("D44").Value = "'1.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.58%  "

$ws.Range("D45").Value = "2.454.52"
$ws.Range("E45").Value = "  -1.12%  "

$ws.Range("D46").Value = "'22.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.58%  "

$ws.Range("E47").Value = "  -2.45%  "

$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("E49").Value = "  -1.46%  "

$ws.Range("D50").Value = "'2.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.09%  "

$ws.Range("E51").Value = "  -2.18%  "
